$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition list) - update "想去人数" (want-to-go count) column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 274
$wsExpo.Range("F6").Value = 63
$wsExpo.Range("F7").Value = 269
$wsExpo.Range("F9").Value = 1996
$wsExpo.Range("F11").Value = 4694

# Sheet "全部类型" (all types combined) - same rows updated
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 274
$wsAll.Range("F8").Value = 63
$wsAll.Range("F9").Value = 269
$wsAll.Range("F13").Value = 1996
$wsAll.Range("F15").Value = 4694
